$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.226199269294739
$ws.Range("B1").Value = 2.81964373588562
$ws.Range("C1").Value = 5.020196437835693
$ws.Range("D1").Value = 2.082340002059937
$ws.Range("E1").Value = 1.160632014274597
